# Generate Report for Handback
#
# For each language sheet (zh-cn, de-de):
#   - Status (col C) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - Two new columns get populated for rows 2/3:
#       F "Latest Target File"   -> mirrors column A (source .md file)
#       G "Latest Handback File" -> mirrors column D (xlf file)
#     both as hyperlinked text using the hyperlink-style font
#     (underline + cornflowerblue, matching the workbook's existing look).
#   - H "Latest Handback DateTime" gets a real timestamp instead of the
#     "0001-01-01 00:00:00" placeholder.

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2        # xlUnderlineStyleSingle
$hyperlinkColor = 15570276     # BGR long for RGB(0x64,0x95,0xED) == FF6495ED

function Set-MirroredHyperlink($ws, $targetCell, $displayText, $url) {
    $ws.Hyperlinks.Add($ws.Range($targetCell), $url, "", "", $displayText) | Out-Null
    $ws.Range($targetCell).Font.Underline = $hyperlinkUnderline
    $ws.Range($targetCell).Font.Color = $hyperlinkColor
}

$sheetsInfo = @(
    @{
        Name = "zh-cn"
        Md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/d31cddde4fd4b82b55677054bc92b9e0500bf0df/e2e/89e5cb83-1b3d-4c8d-b287-01a67560a70f.md"
        Xlf2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4aa11b3e437505c5286986fad48d2e4983a3c0f4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/89e5cb83-1b3d-4c8d-b287-01a67560a70f.4a7fed0aea94a29cded429082adec55f703af684.zh-cn.xlf"
        Xlf2Display = "89e5cb83-1b3d-4c8d-b287-01a67560a70f.4a7fed0aea94a29cded429082adec55f703af684.zh-cn.xlf"
        Md3Url = "https://github.com/OpenLocalizationTest/oltest/blob/d31cddde4fd4b82b55677054bc92b9e0500bf0df/e2e/9a0228e3-9ec4-4a77-b65b-e97f5e96f03e.md"
        Xlf3Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4aa11b3e437505c5286986fad48d2e4983a3c0f4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/9a0228e3-9ec4-4a77-b65b-e97f5e96f03e.647df9de4c9e850f637a00cc17066425878d11db.zh-cn.xlf"
        Xlf3Display = "9a0228e3-9ec4-4a77-b65b-e97f5e96f03e.647df9de4c9e850f637a00cc17066425878d11db.zh-cn.xlf"
        HandbackDateTime = "2016-03-20 04:13:00"
    },
    @{
        Name = "de-de"
        Md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/d31cddde4fd4b82b55677054bc92b9e0500bf0df/e2e/89e5cb83-1b3d-4c8d-b287-01a67560a70f.md"
        Xlf2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cce3dd77298e4d5c4bcb6fe4617950aae9996021/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/89e5cb83-1b3d-4c8d-b287-01a67560a70f.4a7fed0aea94a29cded429082adec55f703af684.de-de.xlf"
        Xlf2Display = "89e5cb83-1b3d-4c8d-b287-01a67560a70f.4a7fed0aea94a29cded429082adec55f703af684.de-de.xlf"
        Md3Url = "https://github.com/OpenLocalizationTest/oltest/blob/d31cddde4fd4b82b55677054bc92b9e0500bf0df/e2e/9a0228e3-9ec4-4a77-b65b-e97f5e96f03e.md"
        Xlf3Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cce3dd77298e4d5c4bcb6fe4617950aae9996021/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/9a0228e3-9ec4-4a77-b65b-e97f5e96f03e.647df9de4c9e850f637a00cc17066425878d11db.de-de.xlf"
        Xlf3Display = "9a0228e3-9ec4-4a77-b65b-e97f5e96f03e.647df9de4c9e850f637a00cc17066425878d11db.de-de.xlf"
        HandbackDateTime = "2016-03-20 04:13:06"
    }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Status: report is now handed back and in sync with en-US.
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Latest Target File (F) / Latest Handback File (G) for row 2,
    # mirroring the source .md (A2) and xlf (D2) hyperlinks.
    Set-MirroredHyperlink $ws "F2" "89e5cb83-1b3d-4c8d-b287-01a67560a70f.md" $info.Md2Url
    Set-MirroredHyperlink $ws "G2" $info.Xlf2Display $info.Xlf2Url

    # Same for row 3.
    Set-MirroredHyperlink $ws "F3" "9a0228e3-9ec4-4a77-b65b-e97f5e96f03e.md" $info.Md3Url
    Set-MirroredHyperlink $ws "G3" $info.Xlf3Display $info.Xlf3Url

    # Latest Handback DateTime (H): replace the 0001-01-01 placeholder.
    $ws.Range("H2").Value = $info.HandbackDateTime
    $ws.Range("H3").Value = $info.HandbackDateTime
}
